# Auto-generated edit script applying the diff's cell-level value changes.
# For each touched row we set the new numeric values for H-N columns as needed,
# and clear any cell that the diff removes entirely (so it drops out of the XML).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 26230.05
$ws.Range("I15").Value = 26230.05
$ws.Range("K15").Value = 78690.14999999999
$ws.Range("M15").Value = -78521.14999999999
$ws.Range("H37").Value = 600
$ws.Range("J37").Value = 600
$ws.Range("L37").Value = 1800
$ws.Range("N37").Value = -2052
$ws.Range("H49").Value = 708.5
$ws.Range("I49").Value = 417
$ws.Range("J49").Value = 1000
$ws.Range("K49").Value = 1251
$ws.Range("L49").Value = 3000
$ws.Range("M49").Value = -1115
$ws.Range("N49").Value = -3272
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H113").Value = 27781210
$ws.Range("I113").Value = 6946850.5
$ws.Range("K113").Value = 6946850.5
$ws.Range("M113").Value = -6943596.5
$ws.Range("H137").Value = 3752.182
$ws.Range("I137").Value = 5254.8
$ws.Range("K137").Value = 15764.4
$ws.Range("M137").Value = -13214.4
$ws.Range("H138").Value = 5705.9824
$ws.Range("I138").Value = 3121.7693
$ws.Range("K138").Value = 9365.3079
$ws.Range("M138").Value = -4225.3079
$ws.Range("H141").Value = 2957.6428
$ws.Range("I141").Value = 1883.125
$ws.Range("K141").Value = 5649.375
$ws.Range("M141").Value = -469.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 224.75
$ws.Range("I4").Value = 224.75
$ws.Range("K4").Value = 224.75
$ws.Range("M4").Value = -108.75
$ws.Range("H61").Value = 7936.758
$ws.Range("I61").Value = 4451.8945
$ws.Range("K61").Value = 4451.8945
$ws.Range("M61").Value = -4239.8945
$ws.Range("H132").Value = 10523.692
$ws.Range("I132").Value = 13422
$ws.Range("K132").Value = 40266
$ws.Range("M132").Value = -37736
$ws.Range("H136").Value = 7936.758
$ws.Range("I136").Value = 4451.8945
$ws.Range("K136").Value = 13355.6835
$ws.Range("M136").Value = -10805.6835

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H28").Value = 48748.25
$ws.Range("J28").Value = 48748.25
$ws.Range("L28").Value = 48748.25
$ws.Range("N28").Value = -49336.25
$ws.Range("H35").Value = 33688.332
$ws.Range("J35").Value = 51000
$ws.Range("L35").Value = 51000
$ws.Range("N35").Value = -51620
$ws.Range("H134").Value = 6491.816
$ws.Range("I134").Value = 3349.739
$ws.Range("K134").Value = 10049.217
$ws.Range("M134").Value = -7514.217000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8442.045
$ws.Range("I31").Value = 4251.846
$ws.Range("J31").Value = 10144.3125
$ws.Range("K31").Value = 4251.846
$ws.Range("L31").Value = 10144.3125
$ws.Range("M31").Value = -3956.846
$ws.Range("N31").Value = -10734.3125
$ws.Range("H34").Value = 8442.045
$ws.Range("I34").Value = 4251.846
$ws.Range("J34").Value = 10144.3125
$ws.Range("K34").Value = 4251.846
$ws.Range("L34").Value = 10144.3125
$ws.Range("M34").Value = -4049.846
$ws.Range("N34").Value = -10548.3125
$ws.Range("H58").Value = 14713312
$ws.Range("I58").Value = 45456840
$ws.Range("K58").Value = 45456840
$ws.Range("M58").Value = -45456637
$ws.Range("H99").Value = 4172.6665
$ws.Range("I99").Value = 2255.5
$ws.Range("K99").Value = 2255.5
$ws.Range("M99").Value = -757.5
$ws.Range("H126").Value = 4172.6665
$ws.Range("I126").Value = 2255.5
$ws.Range("K126").Value = 6766.5
$ws.Range("M126").Value = -4296.5
$ws.Range("H136").Value = 14713312
$ws.Range("I136").Value = 45456840
$ws.Range("K136").Value = 136370520
$ws.Range("M136").Value = -136367970

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I4").Value = 80196290
$ws.Range("J4").Value = 212
$ws.Range("K4").Value = 240588870
$ws.Range("L4").Value = 636
$ws.Range("M4").Value = -240588758
$ws.Range("N4").Value = -860
$ws.Range("H44").Value = 1700
$ws.Range("I44").Value = 100
$ws.Range("J44").Value = 2020
$ws.Range("K44").Value = 300
$ws.Range("L44").Value = 6060
$ws.Range("M44").Value = 98
$ws.Range("N44").Value = -6856
$ws.Range("H87").Value = 62510124
$ws.Range("J87").Value = 11998.923
$ws.Range("L87").Value = 35996.769
$ws.Range("N87").Value = -38492.769
$ws.Range("H90").Value = 62510124
$ws.Range("J90").Value = 11998.923
$ws.Range("L90").Value = 107990.307
$ws.Range("N90").Value = -120470.307
$ws.Range("H122").Value = 1573398.4
$ws.Range("J122").Value = 2118.9
$ws.Range("L122").Value = 19070.1
$ws.Range("N122").Value = -23970.1
$ws.Range("H134").Value = 159005.69
$ws.Range("I134").Value = 159005.69
$ws.Range("K134").Value = 477017.07
$ws.Range("M134").Value = -471947.07
$ws.Range("H137").Value = 184247.73
$ws.Range("J137").Value = 203998.8
$ws.Range("L137").Value = 611996.3999999999
$ws.Range("N137").Value = -622196.3999999999
$ws.Range("H139").Value = 96779.97
$ws.Range("I139").Value = 144712.58
$ws.Range("J139").Value = 5272.273
$ws.Range("K139").Value = 434137.74
$ws.Range("L139").Value = 15816.819
$ws.Range("M139").Value = -428997.74
$ws.Range("N139").Value = -26096.819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1538555.8
$ws.Range("J2").Value = 10000250
$ws.Range("L2").Value = 10000250
$ws.Range("N2").Value = -10000476
$ws.Range("H122").Value = 43843.68
$ws.Range("I122").Value = 61821.94
$ws.Range("K122").Value = 185465.82
$ws.Range("M122").Value = -183015.82
$ws.Range("H132").Value = 4248.033
$ws.Range("I132").Value = 1808.1111
$ws.Range("J132").Value = 7907.9165
$ws.Range("K132").Value = 5424.3333
$ws.Range("L132").Value = 23723.7495
$ws.Range("M132").Value = -2894.3333
$ws.Range("N132").Value = -28783.7495

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3100.8462
$ws.Range("I22").Value = 2060.1667
$ws.Range("J22").Value = 3992.8572
$ws.Range("K22").Value = 2060.1667
$ws.Range("L22").Value = 3992.8572
$ws.Range("M22").Value = -1765.1667
$ws.Range("N22").Value = -4582.8572
$ws.Range("H27").Value = 3100.8462
$ws.Range("I27").Value = 2060.1667
$ws.Range("J27").Value = 3992.8572
$ws.Range("K27").Value = 2060.1667
$ws.Range("L27").Value = 3992.8572
$ws.Range("M27").Value = -1953.1667
$ws.Range("N27").Value = -4206.8572
$ws.Range("H61").Value = 6141.231
$ws.Range("I61").Value = 4750.75
$ws.Range("J61").Value = 6759.222
$ws.Range("K61").Value = 4750.75
$ws.Range("L61").Value = 6759.222
$ws.Range("M61").Value = -4548.75
$ws.Range("N61").Value = -7163.222
$ws.Range("H70").Value = 40997
$ws.Range("J70").Value = 40997
$ws.Range("L70").Value = 40997
$ws.Range("N70").Value = -41537
$ws.Range("H73").Value = 40997
$ws.Range("J73").Value = 40997
$ws.Range("L73").Value = 40997
$ws.Range("N73").Value = -42869
$ws.Range("H113").Value = 6141.231
$ws.Range("I113").Value = 4750.75
$ws.Range("J113").Value = 6759.222
$ws.Range("K113").Value = 4750.75
$ws.Range("L113").Value = 6759.222
$ws.Range("M113").Value = -2580.75
$ws.Range("N113").Value = -11099.222

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 13628.714
$ws.Range("J41").Value = 13628.714
$ws.Range("L41").Value = 13628.714
$ws.Range("N41").Value = -14408.714
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H81").Value = 20007000
$ws.Range("I81").Value = 1250
$ws.Range("K81").Value = 2500
$ws.Range("M81").Value = -1439
$ws.Range("H84").Value = 20007000
$ws.Range("I84").Value = 1250
$ws.Range("K84").Value = 12500
$ws.Range("M84").Value = -7196
$ws.Range("H103").Value = 54797
$ws.Range("J103").Value = 54797
$ws.Range("L103").Value = 54797
$ws.Range("N103").Value = -57141
$ws.Range("H113").Value = 15371.444
$ws.Range("I113").Value = 51086
$ws.Range("K113").Value = 153258
$ws.Range("M113").Value = -151088
$ws.Range("H132").Value = 19256538
$ws.Range("I132").Value = 31258468
$ws.Range("K132").Value = 93775404
$ws.Range("M132").Value = -93772874
$ws.Range("H138").Value = 77995
$ws.Range("J138").Value = 77995
$ws.Range("L138").Value = 77995
$ws.Range("N138").Value = -88275
